$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "United States" column header in F1
$ws.Range("F1").Value = "United States"

# Rows 2-4: full support rows (all countries = 1), add F = 1
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1

# Row 5 (Global tax on millionaires) - updated B:E values, F left empty
$ws.Range("B5").Value = 0.8300395256917
$ws.Range("C5").Value = 0.847826086956522
$ws.Range("D5").Value = 0.879828326180258
$ws.Range("E5").Value = 0.810035842293907

# Row 6 (National tax on millionaires) - updated B:E values, F left empty
$ws.Range("B6").Value = 0.781481481481481
$ws.Range("C6").Value = 0.868421052631579
$ws.Range("D6").Value = 0.892703862660944
$ws.Range("E6").Value = 0.857638888888889

# Row 7 (Payments from high-income countries...)
$ws.Range("B7").Value = 0.686274509803922
$ws.Range("C7").Value = 0.689716312056738
$ws.Range("D7").Value = 0.804147465437788
$ws.Range("E7").Value = 0.667359667359667
$ws.Range("F7").Value = 0.550151975683891

# Row 8 (High-income countries funding renewable energy...)
$ws.Range("B8").Value = 0.803455723542117
$ws.Range("C8").Value = 0.823333333333333
$ws.Range("D8").Value = 0.850678733031674
$ws.Range("E8").Value = 0.798850574712644
$ws.Range("F8").Value = 0.677653902084343

# Row 9 (High-income countries contributing $100 billion...)
$ws.Range("B9").Value = 0.741721854304636
$ws.Range("C9").Value = 0.79020979020979
$ws.Range("D9").Value = 0.775229357798165
$ws.Range("E9").Value = 0.709090909090909
$ws.Range("F9").Value = 0.592648539778449

# Row 10 (Cancellation of low-income countries' public debt)
$ws.Range("B10").Value = 0.528436018957346
$ws.Range("C10").Value = 0.435452793834297
$ws.Range("D10").Value = 0.601965601965602
$ws.Range("E10").Value = 0.587301587301587
$ws.Range("F10").Value = 0.452247191011236

# Row 11 (Democratise international institutions...)
$ws.Range("B11").Value = 0.666666666666667
$ws.Range("C11").Value = 0.695364238410596
$ws.Range("D11").Value = 0.760526315789474
$ws.Range("E11").Value = 0.6953125
$ws.Range("F11").Value = 0.569230769230769

# Row 12 (Removing tariffs on imports from low-income countries)
$ws.Range("B12").Value = 0.595693779904306
$ws.Range("C12").Value = 0.718095238095238
$ws.Range("D12").Value = 0.786885245901639
$ws.Range("E12").Value = 0.844036697247706
$ws.Range("F12").Value = 0.627764127764128

# Row 13 (A minimum wage in all countries...)
$ws.Range("B13").Value = 0.792941176470588
$ws.Range("C13").Value = 0.794117647058824
$ws.Range("D13").Value = 0.807228915662651
$ws.Range("E13").Value = 0.816120906801008
$ws.Range("F13").Value = 0.630372492836676

# Row 14 (Fight tax evasion by creating a global financial register...)
$ws.Range("B14").Value = 0.891304347826087
$ws.Range("C14").Value = 0.847790507364976
$ws.Range("D14").Value = 0.906040268456376
$ws.Range("E14").Value = 0.87378640776699
$ws.Range("F14").Value = 0.6211714132187

# Row 15 (A maximum wealth limit of $10 billion for each human)
$ws.Range("B15").Value = 0.57906976744186
$ws.Range("C15").Value = 0.608247422680412
$ws.Range("D15").Value = 0.618556701030928
$ws.Range("E15").Value = 0.667359667359667
$ws.Range("F15").Value = 0.465809768637532
